# Swap the "Summer Villa" and "Vila Sophia" rows (rows 5 and 6) in the best
# choices table, so that "Vila Sophia" (rating 9,8 / price 750) now appears
# before "Summer Villa" (rating 9,6 / price 828).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (C) holds text values (e.g. "828 ") rather than numbers, so we
# force a text number-format before assigning, otherwise the numeric-looking
# strings would be auto-converted into actual numbers.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C6").NumberFormat = "@"

$ws.Range("A5").Value = "Vila Sophia"
$ws.Range("B5").Value = "9,8"
$ws.Range("C5").Value = "750 "

$ws.Range("A6").Value = "Summer Villa"
$ws.Range("B6").Value = "9,6"
$ws.Range("C6").Value = "828 "

# Restore the original (default) cell formatting on the price cells so the
# only change is the cell content, not its style.
$ws.Range("A1").Copy()
$ws.Range("C5:C6").PasteSpecial(-4122) # xlPasteFormats
